# The deck's slide master currently carries the "Integral" / "Red Violet"
# theme (ppt/theme/theme1.xml). The edit swaps the active design over to the
# stock "Office Theme" palette (the colour scheme that, before this edit,
# only lived unused in ppt/theme/theme2.xml via the notes master).
#
# PowerPoint's Theme object only exposes the 12 theme colour slots for
# in-place editing (ThemeColorScheme.Item(1..12).RGB); the theme/colour
# scheme *names* are read-only through the object model (they only change
# if you save/apply a whole new .thmx), so we drive the colours only.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# Office theme colour scheme, in MsoThemeColorSchemeIndex order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# (values below are the PowerPoint COM RGB() encoding, i.e. R + G*256 + B*65536)

$colors.Item(1).RGB = 0          # dk1      000000
$colors.Item(2).RGB = 16777215   # lt1      FFFFFF
$colors.Item(3).RGB = 6968388    # dk2      44546A
$colors.Item(4).RGB = 15132391   # lt2      E7E6E6
$colors.Item(5).RGB = 13998939   # accent1  5B9BD5
$colors.Item(6).RGB = 3243501    # accent2  ED7D31
$colors.Item(7).RGB = 10855845   # accent3  A5A5A5
$colors.Item(8).RGB = 49407      # accent4  FFC000
$colors.Item(9).RGB = 12874308   # accent5  4472C4
$colors.Item(10).RGB = 4697456   # accent6  70AD47
$colors.Item(11).RGB = 12673797  # hlink    0563C1
$colors.Item(12).RGB = 7491477   # folHlink 954F72
